$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "ArrayList"
$ws.Range("B3").Value = "print(5*8)"
$ws.Range("A4").Value = "Graph"
$ws.Range("B4").Value = "print(15-8)"

$ws.Rows.Item(3).RowHeight = 15.75
$ws.Rows.Item(4).RowHeight = 15.75

$ws.Range("B4").Select()
